$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Max" threshold column (C) values per the new upload.
$ws.Range("C2").Value = 12
$ws.Range("C3").Value = 11
$ws.Range("C4").Value = 1.45
$ws.Range("C5").Value = 25
